# "set sidenav to light mode"
# Adds a third slide, a light-mode (white/bg1 background) copy of slide 1
# (the dark "sidenav"/logo slide), appended to the end of the deck.

$p = $ppt.ActivePresentation

# Slide 1 is the dark-navy (001529) logo slide; slide 2 is the light-gray
# (F2F2F2) variant. Duplicate slide 1 - same logo/group/shapes - then move
# the duplicate to the end of the deck and recolor its background to the
# light "Background 1" theme color (schemeClr val="bg1").
$source = $p.Slides.Item(1)
$source.Duplicate() | Out-Null

# The duplicate is inserted immediately after the source slide (index 2);
# move it to the end so the new sldId lands last in the sldIdLst.
$p.Slides.Item(2).MoveTo($p.Slides.Count)

# Recolor the new last slide's background to the light "Background 1"
# theme color (was the dark navy srgbClr 001529 on slide 1).
$newSlide = $p.Slides.Item($p.Slides.Count)
$newSlide.Background.Fill.ForeColor.SchemeColor = "bg1"
